# Withdraw the key from FanJiaojiao (范娇娇), row 17 of the key log.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# She no longer holds a key - clear the holder cell entirely (not just blank it),
# so the row collapses back to a single populated column like the trailing rows.
$ws.Range("B17").ClearContents()

# Update the "currently X issued / Y in use / Z spare" summary banner:
# one more key went back into the spare pool (16 in use -> 15, 3 spare -> 4).
$summary = $ws.Range("E1")
$summary.Value = "（目前共19把，使用15把，闲置4把）"
$summary.Characters(1, 4).Font.Name = "Noto Sans CJK SC Regular"
$summary.Characters(1, 4).Font.Size = 10
$summary.Characters(5, 2).Font.Name = "Arial"
$summary.Characters(5, 2).Font.Size = 10
$summary.Characters(7, 4).Font.Name = "Noto Sans CJK SC Regular"
$summary.Characters(7, 4).Font.Size = 10
$summary.Characters(11, 2).Font.Name = "Arial"
$summary.Characters(11, 2).Font.Size = 10
$summary.Characters(13, 7).Font.Name = "Noto Sans CJK SC Regular"
$summary.Characters(13, 7).Font.Size = 10

# Leave the cursor resting on the title banner rather than the now-empty B17.
$ws.Range("A1:D1").Select() | Out-Null
